$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("被测试功能3")

# Update rows 9-14: cell value re-shuffle from shared string table changes
$ws.Range("D9").Value = "输入正确的数据类型"
$ws.Range("F10").Value = "弹出升级成功/失败信息"
$ws.Range("C11").Value = "删除用户"
$ws.Range("D11").Value = "选择用户"
$ws.Range("F11").Value = "弹出删除成功/失败信息"
$ws.Range("C12").Value = "指定日期查看统计数据"
$ws.Range("D12").Value = "选择正确日期"
$ws.Range("F12").Value = "显示相关数据"
$ws.Range("C13").Value = "搜索订单"
$ws.Range("D13").Value = "输入正确的数据类型"
$ws.Range("C14").Value = "审核通过订单"
$ws.Range("D14").Value = "订单存在且未通过"
$ws.Range("F14").Value = "弹出审核通过信息"

# Fill in previously-blank rows 15-23 with new test case data
# Row 15
$ws.Range("B15").Value = 7
$ws.Range("C15").Value = "搜索场馆"
$ws.Range("D15").Value = "输入正确的数据类型"
$ws.Range("F15").Value = "显示搜索内容，弹出搜索情况"
# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 7
$ws.Range("C16").Value = "添加场馆"
$ws.Range("D16").Value = "输入正确的数据格式"
$ws.Range("F16").Value = "弹出添加成功/失败信息"
# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 7
$ws.Range("C17").Value = "编辑场馆"
$ws.Range("D17").Value = "输入正确的数据格式"
$ws.Range("F17").Value = "弹出编辑成功/失败信息"
# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 8
$ws.Range("C18").Value = "添加新闻"
$ws.Range("D18").Value = "输入正确的数据格式"
$ws.Range("F18").Value = "弹出添加成功/失败信息"
# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 8
$ws.Range("C19").Value = "删除新闻"
$ws.Range("D19").Value = "选择新闻不为空"
$ws.Range("F19").Value = "弹出删除成功/失败信息"
$ws.Range("G19").Value = "是"
# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 8
$ws.Range("C20").Value = "编辑新闻"
$ws.Range("D20").Value = "输入正确的数据格式"
$ws.Range("F20").Value = "弹出编辑成功/失败信息"
$ws.Range("G20").Value = "是"
# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 8
$ws.Range("C21").Value = "查看图片"
$ws.Range("D21").Value = "新闻存在图片"
$ws.Range("F21").Value = "弹出显示图片的对话框"
$ws.Range("G21").Value = "是"
# Row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 9
$ws.Range("C22").Value = "搜索留言"
$ws.Range("D22").Value = "输入正确的数据类型"
$ws.Range("F22").Value = "显示搜索内容，弹出搜索情况"
$ws.Range("G22").Value = "是"
# Row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 9
$ws.Range("C23").Value = "审核通过留言"
$ws.Range("D23").Value = "留言存在且未通过"
$ws.Range("F23").Value = "弹出审核通过信息"
$ws.Range("G23").Value = "是"

# Update view: selection moved to A23 (scroll position A2 not representable via COM)
$ws.Activate()
$ws.Range("A23").Select()
